# Metro Divisions - start looking at LA (Local Area) database.
# Flags the 11 "Metropolitan Statistical Area" header rows with a new
# highlight color and a marker "1" in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-B cell is a top-level "... Metropolitan Statistical Area"
# / "... Metropolitan NECTA" header (as opposed to a Metropolitan Division
# detail row).
$headerRows = @(3, 9, 13, 17, 21, 26, 32, 38, 43, 47, 51)

# 1) Mark each header row with a 1 in column C.
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 3).Value = 1
}

# 2) Give the header rows' B cell a new highlight fill (accent2/orange family,
#    matching the existing accent6/green highlight already used on the
#    Metropolitan Division detail rows).
$firstHeader = $ws.Range("B3")
$firstHeader.Interior.ThemeColor = 6
$firstHeader.Copy()
foreach ($r in $headerRows) {
    if ($r -ne 3) {
        $ws.Cells.Item($r, 2).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    }
}

# 3) B60 picks up the same highlight style already used on the other
#    Metropolitan/NECTA Division detail rows (e.g. B4).
$ws.Range("B4").Copy()
$ws.Range("B60").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# 4) Update the view: select B3:B61 (started reviewing the LA / Local Area
#    rows) instead of the previous scrolled-down single-cell selection.
$ws.Range("B3:B61").Select()
